$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44883
$ws.Range("J2").Value = 290
$ws.Range("K2").Value = 1400
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = 1434
$ws.Range("P2").Value = 1434

$ws.Range("D3").Value = 44895
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 1200
$ws.Range("L3").Value = 1300
$ws.Range("M3").Value = 1255
$ws.Range("P3").Value = 1255

$ws.Range("D4").Value = 44893
$ws.Range("J4").Value = 3300
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 1300
$ws.Range("M4").Value = 1261
$ws.Range("P4").Value = 1261

$ws.Range("D5").Value = 44638
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 2800
$ws.Range("M5").Value = 2650
$ws.Range("P5").Value = 2650

$ws.Range("D6").Value = 44175
$ws.Range("J6").Value = 1400
$ws.Range("K6").Value = 1900
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 1950
$ws.Range("P6").Value = 1950

$ws.Range("D7").Value = 44210
$ws.Range("J7").Value = 1450
$ws.Range("K7").Value = 1600
$ws.Range("L7").Value = 1700
$ws.Range("M7").Value = 1650
$ws.Range("P7").Value = 1650

$ws.Range("D8").Value = 44200
$ws.Range("J8").Value = 1500
$ws.Range("M8").Value = 1450
$ws.Range("P8").Value = 1450

$ws.Range("D9").Value = 44907
$ws.Range("J9").Value = 2300
$ws.Range("K9").Value = 900
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 952
$ws.Range("P9").Value = 952

$ws.Range("D10").Value = 44537
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 1300
$ws.Range("L10").Value = 1400
$ws.Range("M10").Value = 1350
$ws.Range("P10").Value = 1350
